# Update the Senegal Pre-TAS FTS form from V2 to V3.
#
# Content changes:
#   - survey!D6  (hint::French for d_cluster_id): "... chaque école" -> "... chaque village"
#   - survey!B9  (name of the "begin repeat" row): sn_lf_f_2407_2 -> sn_lf_f_2407_3
#   - settings!A2 (form_title): "...FTS V2" -> "...FTS V3"
#   - settings!B2 (form_id): sn_lf_pretas_20407_2_fts_v2 -> sn_lf_pretas_20407_2_fts_v3

$wb = $excel.ActiveWorkbook

$wsSurvey   = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")

# --- content edits -------------------------------------------------------

$wsSurvey.Range("D6").Value = "Le code a trois chiffre associé à chaque village"
$wsSurvey.Range("B9").Value = "sn_lf_f_2407_3"

$wsSettings.Range("A2").Value = "(2024 Juillet) 2. Pre-TAS - Formulaire Résultat FTS V3"
$wsSettings.Range("B2").Value = "sn_lf_pretas_20407_2_fts_v3"

# --- selection / active sheet --------------------------------------------
# Final state: "survey" is the active tab with D6 selected; "settings" is
# left with A2 selected (but not active). Select settings!A2 first, then
# finish on survey!D6 so survey ends up the active sheet.

$wsSettings.Range("A2").Select()
$wsSurvey.Activate()
$wsSurvey.Range("D6").Select()
